# ERK main axis colored
# Slide 3 ("schemes") has a row of 7 small legend rectangles that were all
# using a themed accent6 color at varying alpha levels. Each one gets its
# own explicit RGB color (still with an alpha/transparency), and the
# caption below the row changes from "velocity axis direction" to
# "main axis direction".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# Shape index (in the slide's Shapes collection) -> (new RGB hex, new alpha in 1000ths of a percent)
# alpha 50000 = 50% opaque = 50% transparent (Fill.Transparency = 0.5)
# alpha 75000 = 75% opaque = 25% transparent (Fill.Transparency = 0.25)
$changes = @(
    @{ Index = 2;  Rgb = 0x4472C4; Alpha = 50000 },
    @{ Index = 4;  Rgb = 0x235888; Alpha = 50000 },
    @{ Index = 6;  Rgb = 0x2C70AE; Alpha = 75000 },
    @{ Index = 8;  Rgb = 0xBA7741; Alpha = 50000 },
    @{ Index = 10; Rgb = 0xE08F4E; Alpha = 50000 },
    @{ Index = 12; Rgb = 0xFFB858; Alpha = 75000 },
    @{ Index = 14; Rgb = 0xFFC000; Alpha = 50000 }
)

foreach ($ch in $changes) {
    $sh = $s.Shapes.Item($ch.Index)

    $hex = $ch.Rgb
    $r = [math]::Floor($hex / 0x10000) -band 0xFF
    $g = [math]::Floor($hex / 0x100) -band 0xFF
    $b = $hex -band 0xFF

    # COM ForeColor.RGB takes a BGR-packed integer (same as the VBA RGB() fn).
    $bgr = ($b * 65536) + ($g * 256) + $r

    $sh.Fill.ForeColor.RGB = $bgr
    $sh.Fill.Transparency = (100000 - $ch.Alpha) / 100000
}

# Caption text: "velocity axis direction" -> "main axis direction"
$caption = $s.Shapes.Item(17)
$textRange = $caption.TextFrame.TextRange
$velocity = $textRange.Characters(1, 8)
$velocity.Text = "main"
